$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Errol"
$ws.Range("B2").Value = "Lueilwitz"
$ws.Range("C2").Value = "fatimah.kling@gmail.com"
$ws.Range("D2").Value = "ffy04rn7msl"
